$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial value that was bumped by one day
# (45607 -> 45608, i.e. 2024-11-11 -> 2024-11-12) for every data row (2..32).
for ($r = 2; $r -le 32; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45607) {
        $cell.Value2 = 45608
    }
}
